# Adds two new instrument models/instruments/deployments:
#   - Roche Modular P Chemistry Analyzer
#   - Tosoh G7 Automated HPLC Analyzer
# to the InstrumentModels, Instruments and Deployments sheets.

$wb = $excel.ActiveWorkbook

$wsModels = $wb.Worksheets.Item("InstrumentModels")
$wsInstr  = $wb.Worksheets.Item("Instruments")
$wsDeploy = $wb.Worksheets.Item("Deployments")

# =================================================================
# Roche Modular P Chemistry Analyzer
# =================================================================

# --- InstrumentModels: insert new row 4, copying the formatting of row 3
# (an existing fully-styled model row).
$wsModels.Rows.Item(3).Copy()
$wsModels.Rows.Item(4).Insert()
$wsModels.Cells.Item(4,1).Value = "nhanes-kb:ROCHE-MODULAR-P-CHEMISTRY-ANALYZER"
$wsModels.Cells.Item(4,2).Value = "vstoi:PhysicalInstrument"
$wsModels.Cells.Item(4,3).Value = "Roche Modular P Chemistry Analyzer"
$wsModels.Cells.Item(4,4).Value = "Roche Diagnostics"

# --- Instruments: insert new row 73, copying the formatting of row 72 (an
# existing fully-styled instrument row).
$wsInstr.Rows.Item(72).Copy()
$wsInstr.Rows.Item(73).Insert()
$wsInstr.Cells.Item(73,1).Value = "nhanes-kb:INS-ROCHE-MODULAR-P-CHEMISTRY-ANALYZER"
$wsInstr.Cells.Item(73,2).Value = "nhanes-kb:ROCHE-MODULAR-P-CHEMISTRY-ANALYZER"
$wsInstr.Cells.Item(73,3).Value = "Generic Roche Modular P Chemistry Analyzer"

# --- Deployments: insert new row 28, copying the formatting of row 27 (an
# existing fully-styled deployment row); only the hasURI and hasInstrument
# columns differ between deployment rows.
$wsDeploy.Rows.Item(27).Copy()
$wsDeploy.Rows.Item(28).Insert()
$wsDeploy.Cells.Item(28,1).Value = "nhanes-kb:DPL-ROCHE-MODULAR-P-CHEMISTRY-ANALYZER"
$wsDeploy.Cells.Item(28,4).Value = "nhanes-kb:INS-ROCHE-MODULAR-P-CHEMISTRY-ANALYZER"
$wsDeploy.Rows.Item(28).RowHeight = 15.75

# =================================================================
# Tosoh G7 Automated HPLC Analyzer
# =================================================================

# --- InstrumentModels: insert new row 6, copying the formatting of row 6
# (after the previous insert this is the shifted-down "Tosoh G8
# Glycohemoglobin Analyzer" row, which carries the correct style).
$wsModels.Rows.Item(6).Copy()
$wsModels.Rows.Item(6).Insert()
$wsModels.Cells.Item(6,1).Value = "nhanes-kb:TOSOH-G7-AUTOMATED-HPLC-ANALYZER"
$wsModels.Cells.Item(6,2).Value = "vstoi:PhysicalInstrument"
$wsModels.Cells.Item(6,3).Value = "Tosoh G7 Automated HPLC Analyzer"
$wsModels.Cells.Item(6,4).Value = "Tosoh Medics, Inc."

# --- Instruments: insert new row 74, copying the formatting of row 74
# (after the previous insert this is the shifted-down "Tosoh G8" instrument
# row, which carries the correct style).
$wsInstr.Rows.Item(74).Copy()
$wsInstr.Rows.Item(74).Insert()
$wsInstr.Cells.Item(74,1).Value = "nhanes-kb:INS-TOSOH-G7-AUTOMATED-HPLC-ANALYZER"
$wsInstr.Cells.Item(74,2).Value = "nhanes-kb:TOSOH-G7-AUTOMATED-HPLC-ANALYZER"
$wsInstr.Cells.Item(74,3).Value = "Generic Tosoh G7 Automated HPLC Analyzer"

# --- Deployments: insert new row 29, copying the formatting of row 29
# (after the previous insert this is the shifted-down "Tosoh G8" deployment
# row, which carries the correct style).
$wsDeploy.Rows.Item(29).Copy()
$wsDeploy.Rows.Item(29).Insert()
$wsDeploy.Cells.Item(29,1).Value = "nhanes-kb:DPL-TOSOH-G7-AUTOMATED-HPLC-ANALYZER"
$wsDeploy.Cells.Item(29,4).Value = "nhanes-kb:INS-TOSOH-G7-AUTOMATED-HPLC-ANALYZER"
$wsDeploy.Rows.Item(29).RowHeight = 15.75
